$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an always-default-styled cell (never written to) so we can
# force text storage for numeric-looking strings (NumberFormat "@")
# while restoring the resulting cell style to match the original
# (unstyled) data cells once the value has been written.
$defaultStyle = $ws.Range("A1").Style

# Rows 18 and 19: Litecoin and ShibaInu swap places, with updated price/volume values
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001053"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("E18").Style = $defaultStyle

$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.54"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("E19").Style = $defaultStyle

# Update Price (D) and Volume(1h) (E) columns for remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.024.16"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.897.36"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E3").Style = $defaultStyle
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.40%  "
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.06"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E5").Style = $defaultStyle
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4708"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E7").Style = $defaultStyle
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3943"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E8").Style = $defaultStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.67"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E9").Style = $defaultStyle
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08086"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E10").Style = $defaultStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.026"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E11").Style = $defaultStyle
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.95"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.894.06"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.996"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E14").Style = $defaultStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.161"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("E15").Style = $defaultStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.018"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E16").Style = $defaultStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06810"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.29"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("E20").Style = $defaultStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.013"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E21").Style = $defaultStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.027.04"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("E22").Style = $defaultStyle
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.535"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E23").Style = $defaultStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E24").Style = $defaultStyle
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.350"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E25").Style = $defaultStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.091.60"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E26").Style = $defaultStyle
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.89"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("E27").Style = $defaultStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.13"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E28").Style = $defaultStyle
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.100"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E29").Style = $defaultStyle
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.501"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("E30").Style = $defaultStyle
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.19"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E31").Style = $defaultStyle
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9774"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("E32").Style = $defaultStyle
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09535"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E33").Style = $defaultStyle
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.647"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E34").Style = $defaultStyle
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.415"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("E35").Style = $defaultStyle
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.392"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E36").Style = $defaultStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06161"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("E37").Style = $defaultStyle
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02267"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E38").Style = $defaultStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.222"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E39").Style = $defaultStyle
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.119"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E40").Style = $defaultStyle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6019"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E41").Style = $defaultStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1897"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E42").Style = $defaultStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.34"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.267"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("E44").Style = $defaultStyle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5732"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E45").Style = $defaultStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("E46").Style = $defaultStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.413"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E47").Style = $defaultStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.948"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E48").Style = $defaultStyle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06945"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("E49").Style = $defaultStyle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.30"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("E50").Style = $defaultStyle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.074"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.82%  "
$ws.Range("E51").Style = $defaultStyle
